$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet 1 ("Project 1") ----
# Update header row with the new "Team/Email/Phone/Budget" labels
$ws1.Range("A1").Value = "Team B P1"
$ws1.Range("B1").Value = "Email B P1"
$ws1.Range("C1").Value = "Phone B P1"
$ws1.Range("D1").Value = "Budget B P1"

# Remove the old second data row entirely
$ws1.Range("A2:C2").ClearContents()

# Big header styling
$ws1.Range("A1:D1").Font.Size = 24
$ws1.Rows.Item(1).RowHeight = 31.5

# Column widths (best fit for the new header text)
$ws1.Columns.Item(1).ColumnWidth = 20.8776
$ws1.Columns.Item(2).ColumnWidth = 20.59245
$ws1.Columns.Item(3).ColumnWidth = 22.59245
$ws1.Columns.Item(4).ColumnWidth = 23.8776

$ws1.Range("A1:D1").Select()

# ---- Sheet 2 ("Project 2") ----
# Remove the old first data row entirely
$ws2.Range("A1:C1").ClearContents()

# Update second row with the new "Team/Email/Phone/Budget" labels
$ws2.Range("A2").Value = "Team B P2"
$ws2.Range("B2").Value = "Email B P2"
$ws2.Range("C2").Value = "Phone B P2"
$ws2.Range("D2").Value = "Budget B P2"

# Header styling
$ws2.Range("A2:D2").Font.Size = 18
$ws2.Rows.Item(2).RowHeight = 23.25

# Column widths (best fit for the new header text)
$ws2.Columns.Item(1).ColumnWidth = 15.30729
$ws2.Columns.Item(2).ColumnWidth = 15.30729
$ws2.Columns.Item(3).ColumnWidth = 16.45182
$ws2.Columns.Item(4).ColumnWidth = 17.59245

# Make "Project 2" the active tab/sheet with A2:D2 selected
$ws2.Activate()
$ws2.Range("A2:D2").Select()
